# BINGO_cc.xlsx — "Added a cliche and tweaked others"
#
# Changes applied:
#  - "list" sheet (column A, the master list of cliches/phrases): five
#    entries are reworded. Writing them in this specific order reproduces
#    the same shared-string table ordering as the authoritative diff
#    (new strings appended as: "Can you / everyone see my screen?",
#    "I was multi-tasking", "Echo / feedback / cutting out",
#    "Wind / road noise", "Reach out to  _____"), while the two obsolete
#    phrasings ("Can you repeat? I was multi-tasking" and
#    "Can you reach out to  _____?") fall out of the shared-string table
#    since nothing references them anymore.
#  - The sheet's active selection moves from A45 to A28 (the view was
#    scrolled back up a bit after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

$ws.Range("A22").Value = "Can you / everyone see my screen?"
$ws.Range("A31").Value = "I was multi-tasking"
$ws.Range("A10").Value = "Echo / feedback / cutting out"
$ws.Range("A19").Value = "Wind / road noise"
$ws.Range("A27").Value = "Reach out to  _____"

# Scroll/selection: the saved view now shows row 29 at the top with A28
# selected (previously topLeftCell A25 / selection A45).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A28").Select()
